$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1756658.6
$ws.Range("J17").Value = 1788015.8
$ws.Range("L17").Value = 5364047.4
$ws.Range("N17").Value = -5364383.4

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 252.71428
$ws.Range("I39").Value = 166.33333
$ws.Range("J39").Value = 317.5
$ws.Range("K39").Value = 498.99999
$ws.Range("L39").Value = 952.5
$ws.Range("M39").Value = -202.99999
$ws.Range("N39").Value = -1544.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2962
$ws.Range("I62").Value = 2400.3635
$ws.Range("J62").Value = 3991.6667
$ws.Range("K62").Value = 2400.3635
$ws.Range("L62").Value = 3991.6667
$ws.Range("M62").Value = -1776.3635
$ws.Range("N62").Value = -5239.6667

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 2962
$ws.Range("I65").Value = 2400.3635
$ws.Range("J65").Value = 3991.6667
$ws.Range("K65").Value = 12001.8175
$ws.Range("L65").Value = 19958.3335
$ws.Range("M65").Value = -8881.817499999999
$ws.Range("N65").Value = -26198.3335

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 3487.0667
$ws.Range("I116").Value = 1465.8334
$ws.Range("J116").Value = 4834.5557
$ws.Range("K116").Value = 1465.8334
$ws.Range("L116").Value = 4834.5557
$ws.Range("M116").Value = 1976.1666
$ws.Range("N116").Value = -11718.5557

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 1026.4
$ws.Range("I125").Value = 1044
$ws.Range("K125").Value = 9396
$ws.Range("M125").Value = -6936

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 3251.5
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 3251.5
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 9754.5
$ws.Range("N132").Value = -14814.5
$ws.Range("M132").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 47851.727
$ws.Range("I137").Value = 2549.1333
$ws.Range("J137").Value = 144928.72
$ws.Range("K137").Value = 7647.3999
$ws.Range("L137").Value = 434786.16
$ws.Range("M137").Value = -5097.3999
$ws.Range("N137").Value = -439886.16

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3582.9333
$ws.Range("I61").Value = 2984.4
$ws.Range("K61").Value = 2984.4
$ws.Range("M61").Value = -2772.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 500000700
$ws.Range("I74").Value = 1000000000
$ws.Range("J74").Value = 1414
$ws.Range("K74").Value = 1000000000
$ws.Range("L74").Value = 1414
$ws.Range("M74").Value = -999999126
$ws.Range("N74").Value = -3162

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 500000700
$ws.Range("I77").Value = 1000000000
$ws.Range("J77").Value = 1414
$ws.Range("K77").Value = 5000000000
$ws.Range("L77").Value = 7070
$ws.Range("M77").Value = -4999995632
$ws.Range("N77").Value = -15806

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 909.0454999999999
$ws.Range("I97").Value = 815.7368
$ws.Range("J97").Value = 1500
$ws.Range("K97").Value = 815.7368
$ws.Range("L97").Value = 1500
$ws.Range("M97").Value = -319.7368
$ws.Range("N97").Value = -2492

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 13326.613
$ws.Range("I132").Value = 1597.909
$ws.Range("J132").Value = 48512.727
$ws.Range("K132").Value = 4793.727000000001
$ws.Range("L132").Value = 145538.181
$ws.Range("M132").Value = -2263.727000000001
$ws.Range("N132").Value = -150598.181

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3582.9333
$ws.Range("I136").Value = 2984.4
$ws.Range("K136").Value = 8953.200000000001
$ws.Range("M136").Value = -6403.200000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 48592.87
$ws.Range("I134").Value = 58296.633
$ws.Range("J134").Value = 2500
$ws.Range("K134").Value = 174889.899
$ws.Range("L134").Value = 7500
$ws.Range("M134").Value = -172354.899
$ws.Range("N134").Value = -12570

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8405.710999999999
$ws.Range("I31").Value = 10626.032
$ws.Range("J31").Value = 3489.2856
$ws.Range("K31").Value = 10626.032
$ws.Range("L31").Value = 3489.2856
$ws.Range("M31").Value = -10331.032
$ws.Range("N31").Value = -4079.2856

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 8405.710999999999
$ws.Range("I34").Value = 10626.032
$ws.Range("J34").Value = 3489.2856
$ws.Range("K34").Value = 10626.032
$ws.Range("L34").Value = 3489.2856
$ws.Range("M34").Value = -10424.032
$ws.Range("N34").Value = -3893.2856

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 16276.546
$ws.Range("I58").Value = 1197.0454
$ws.Range("K58").Value = 1197.0454
$ws.Range("M58").Value = -994.0454

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 2454.6667
$ws.Range("J94").Value = 3768.6667
$ws.Range("L94").Value = 3768.6667
$ws.Range("N94").Value = -4670.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 19246.232
$ws.Range("I132").Value = 26556.8
$ws.Range("K132").Value = 79670.39999999999
$ws.Range("M132").Value = -77140.39999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1377
$ws.Range("I134").Value = 1032.0667
$ws.Range("K134").Value = 3096.2001
$ws.Range("M134").Value = -561.2001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 16276.546
$ws.Range("I136").Value = 1197.0454
$ws.Range("K136").Value = 3591.1362
$ws.Range("M136").Value = -1041.1362

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1330.7142
$ws.Range("J5").Value = 1668.3334
$ws.Range("L5").Value = 5005.0002
$ws.Range("N5").Value = -5229.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 125000570
$ws.Range("I86").Value = 975
$ws.Range("J86").Value = 250000160
$ws.Range("K86").Value = 2925
$ws.Range("L86").Value = 750000480
$ws.Range("M86").Value = -1739
$ws.Range("N86").Value = -750002852

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H89").Value = 125000570
$ws.Range("I89").Value = 975
$ws.Range("J89").Value = 250000160
$ws.Range("K89").Value = 8775
$ws.Range("L89").Value = 2250001440
$ws.Range("M89").Value = -2847
$ws.Range("N89").Value = -2250013296

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H116").Value = 999.6667
$ws.Range("I116").Value = 999.6667
$ws.Range("K116").Value = 2999.0001
$ws.Range("M116").Value = 442.9998999999998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 709.5
$ws.Range("J122").Value = 828.3
$ws.Range("L122").Value = 7454.7
$ws.Range("N122").Value = -12354.7

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 736.28864
$ws.Range("J131").Value = 749.68134
$ws.Range("L131").Value = 2249.04402
$ws.Range("N131").Value = -12329.04402

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 1595.6177
$ws.Range("I134").Value = 1359.3103
$ws.Range("J134").Value = 2966.2
$ws.Range("K134").Value = 4077.9309
$ws.Range("L134").Value = 8898.599999999999
$ws.Range("M134").Value = 992.0690999999997
$ws.Range("N134").Value = -19038.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 1330.7142
$ws.Range("J135").Value = 1668.3334
$ws.Range("L135").Value = 15015.0006
$ws.Range("N135").Value = -20085.0006

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 4116.8823
$ws.Range("I126").Value = 3217.8096
$ws.Range("J126").Value = 5569.231
$ws.Range("K126").Value = 9653.4288
$ws.Range("L126").Value = 16707.693
$ws.Range("M126").Value = -7183.4288
$ws.Range("N126").Value = -21647.693

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 100015.56
$ws.Range("I132").Value = 107723.7
$ws.Range("K132").Value = 323171.1
$ws.Range("M132").Value = -320641.1

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1970
$ws.Range("I22").Value = 1468.1818
$ws.Range("J22").Value = 3350
$ws.Range("K22").Value = 1468.1818
$ws.Range("L22").Value = 3350
$ws.Range("M22").Value = -1173.1818
$ws.Range("N22").Value = -3940

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 1970
$ws.Range("I27").Value = 1468.1818
$ws.Range("J27").Value = 3350
$ws.Range("K27").Value = 1468.1818
$ws.Range("L27").Value = 3350
$ws.Range("M27").Value = -1361.1818
$ws.Range("N27").Value = -3564

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3591.4167
$ws.Range("I132").Value = 2442.7144
$ws.Range("J132").Value = 5199.6
$ws.Range("K132").Value = 7328.1432
$ws.Range("L132").Value = 15598.8
$ws.Range("M132").Value = -4798.1432
$ws.Range("N132").Value = -20658.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 33312.562
$ws.Range("I136").Value = 63750.125
$ws.Range("J136").Value = 2875
$ws.Range("K136").Value = 191250.375
$ws.Range("L136").Value = 8625
$ws.Range("M136").Value = -188700.375
$ws.Range("N136").Value = -13725

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 4492
$ws.Range("I14").Value = 4004
$ws.Range("K14").Value = 4004
$ws.Range("M14").Value = -3836

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H119").Value = 24500
$ws.Range("J119").Value = 24500
$ws.Range("L119").Value = 24500
$ws.Range("N119").Value = -34176

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 974.6842
$ws.Range("I126").Value = 833.9
$ws.Range("J126").Value = 1131.1111
$ws.Range("K126").Value = 2501.7
$ws.Range("L126").Value = 3393.3333
$ws.Range("M126").Value = -31.69999999999982
$ws.Range("N126").Value = -8333.3333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2155.9167
$ws.Range("I132").Value = 1247.25
$ws.Range("K132").Value = 3741.75
$ws.Range("M132").Value = -1211.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 45457188
$ws.Range("I136").Value = 100003010
$ws.Range("J136").Value = 2342.4167
$ws.Range("K136").Value = 300009030
$ws.Range("L136").Value = 7027.250100000001
$ws.Range("M136").Value = -300006480
